# Generate Report for Handoff
# The "4468ce90-7ace-4bbe-927a-4e931ed44153" entity has been handed off and is no
# longer pending handback, so its row is removed from every sheet. The remaining
# "43ad64a5-f232-48e2-909d-d15f4e532309" entity's status moves from "Handed back:
# in sync with en-US" to "Ready for handoff", and its handoff timestamps on the
# zh-cn / de-de sheets are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Update the status text for the remaining entity.
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"

# Drop the row belonging to the handed-off entity; the .localization-config
# row below it shifts up to row 3.
$ov.Rows.Item(3).Delete()

# Hyperlinks are not renumbered automatically when rows move, so rebuild them.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3e3b86e2c9d4f67e5c8b9ae844c27aa3a5414bf9/e2e/43ad64a5-f232-48e2-909d-d15f4e532309.md", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3e3b86e2c9d4f67e5c8b9ae844c27aa3a5414bf9/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = "Ready for handoff"
$zh.Range("D2").Value = "2016-03-10 18:49:23"

$zh.Rows.Item(3).Delete()

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3e3b86e2c9d4f67e5c8b9ae844c27aa3a5414bf9/e2e/43ad64a5-f232-48e2-909d-d15f4e532309.md", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.md")
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd26983a2cfee9d0808361ca4b9705537ff39abb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.zh-cn.xlf", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/1b61c74dae13f220096d628cfff56055e16decbe/e2e/43ad64a5-f232-48e2-909d-d15f4e532309.md", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.md")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e06654ad13a48316e06471f03ed348daa7e02f16/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.zh-cn.xlf", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3e3b86e2c9d4f67e5c8b9ae844c27aa3a5414bf9/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = "Ready for handoff"
$de.Range("D2").Value = "2016-03-10 18:49:28"

$de.Rows.Item(3).Delete()

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3e3b86e2c9d4f67e5c8b9ae844c27aa3a5414bf9/e2e/43ad64a5-f232-48e2-909d-d15f4e532309.md", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.md")
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8c21e0a0f7bba3a842fdadc383e127ca1d1c0142/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.de-de.xlf", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/407592cdd23b932c8112cf72d418cef6d115f1c0/e2e/43ad64a5-f232-48e2-909d-d15f4e532309.md", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.md")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/486c6be6bce9c98128b7973c2325475e6e7a6471/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.de-de.xlf", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3e3b86e2c9d4f67e5c8b9ae844c27aa3a5414bf9/.localization-config", "", "", ".localization-config")
